$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = "new price text"; E = "new volume text" }
$updates = @{
    2  = @{ D = "63.624.17";  E = "  +5.25%  " }
    3  = @{ D = "3.063.65";   E = "  +3.77%  " }
    4  = @{ E = "  +0.14%  " }
    5  = @{ D = "550.27";     E = "  +6.07%  " }
    6  = @{ D = "138.96";     E = "  +7.60%  " }
    8  = @{ D = "3.058.68";   E = "  +3.69%  " }
    9  = @{ E = "  +4.69%  " }
    10 = @{ D = "0.150";      E = "  +3.06%  " }
    11 = @{ D = "6.23";       E = "  +1.54%  " }
    12 = @{ D = "0.454";      E = "  +4.72%  " }
    13 = @{ E = "  +4.70%  " }
    14 = @{ D = "34.90";      E = "  +6.51%  " }
    15 = @{ D = "3.565.43";   E = "  +3.86%  " }
    16 = @{ D = "63.667.43";  E = "  +5.44%  " }
    17 = @{ D = "3.067.37";   E = "  +4.09%  " }
    18 = @{ E = "  -0.34%  " }
    19 = @{ D = "6.74";       E = "  +5.01%  " }
    20 = @{ D = "483.78";     E = "  +6.54%  " }
    21 = @{ D = "13.55";      E = "  +4.78%  " }
    22 = @{ E = "  +3.12%  " }
    23 = @{ E = "  +6.95%  " }
    24 = @{ D = "81.59";      E = "  +5.21%  " }
    25 = @{ E = "  +8.61%  " }
    26 = @{ E = "  -0.07%  " }
    27 = @{ D = "2.75";       E = "  +5.98%  " }
    28 = @{ D = "7.98";       E = "  +5.05%  " }
    29 = @{ E = "  +10.21%  " }
    30 = @{ D = "0.999";      E = "  +0.17%  " }
    31 = @{ E = "  +4.09%  " }
    32 = @{ D = "1.15";       E = "  +2.84%  " }
    33 = @{ E = "  +9.41%  " }
    34 = @{ D = "5.75";       E = "  +9.28%  " }
    35 = @{ D = "55.28";      E = "  +1.30%  " }
    36 = @{ D = "5.99";       E = "  +5.30%  " }
    37 = @{ D = "466.53";     E = "  +4.81%  " }
    38 = @{ D = "3.159.71";   E = "  +0.62%  " }
    39 = @{ E = "  +6.20%  " }
    40 = @{ D = "0.0395";     E = "  +5.09%  " }
    41 = @{ D = "0.119";      E = "  +3.31%  " }
    42 = @{ E = "  +4.22%  " }
    43 = @{ D = "28.49";      E = "  +14.49%  " }
    44 = @{ D = "2.53";       E = "  +5.88%  " }
    45 = @{ D = "0.251";      E = "  +4.44%  " }
    47 = @{ E = "  +7.04%  " }
    48 = @{ E = "  +2.99%  " }
    49 = @{ D = "0.0₃0510";   E = "  +2.27%  " }
    50 = @{ D = "116.40";     E = "  -2.15%  " }
    51 = @{ E = "  +6.54%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]

    if ($rowData.ContainsKey("D")) {
        # Price column holds text like "63.624.17" or "138.96" that must stay
        # text rather than being auto-parsed as a number/date by Excel.
        # Force the text format, assign, then restore the default style so
        # no stray formatting is left behind on the cell.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["D"]
        $cell.Style = "Normal"
    }

    if ($rowData.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $rowData["E"]
    }
}
